# Generate Report for Handback
# Updates the localization-status report after a successful handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File / Latest Handback DateTime are filled in
#    for the zh-cn and de-de language rows, with a hyperlink on the target file.
#  - A few columns are widened so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFile = "940c231c-cdc6-458a-8955-ab89043e0388.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/4efc6ebe575d02e7076ac31591956137185f3a43/e2e/940c231c-cdc6-458a-8955-ab89043e0388.md"

# ---------------------------------------------------------------------------
# Overview sheet: the Status summary columns for zh-cn / de-de (E2, F2)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("J2").Value = "940c231c-cdc6-458a-8955-ab89043e0388.989096b01a9e7fe8fde66d198afa168a9d196fcb.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-14 17:18:42"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, "", "", $targetFile)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("J2").Value = "940c231c-cdc6-458a-8955-ab89043e0388.989096b01a9e7fe8fde66d198afa168a9d196fcb.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-14 17:18:53"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, "", "", $targetFile)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
